$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: "Save" - copy formatting from G1 (bold/border/center header style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for the new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
